# Generate Report for Handoff
#
# The localization status report moves from "In Translation" to
# "Ready for handoff" and the associated timestamps are refreshed to
# reflect the moment the handoff package was generated.
#
#   Overview!E2 (zh-cn status)          : In Translation -> Ready for handoff
#   Overview!F2 (de-de status)          : In Translation -> Ready for handoff
#   Overview!G2 (Latest HO Xliff date)  : 2016-08-22 02:49:20 -> 2016-08-22 02:49:57
#   zh-cn!C2    (Status)                : In Translation -> Ready for handoff
#   zh-cn!H2    (Latest Handoff Datetime): 2016-08-22 02:49:15 -> 2016-08-22 02:49:53
#   de-de!C2    (Status)                : In Translation -> Ready for handoff
#   de-de!H2    (Latest Handoff Datetime): 2016-08-22 02:49:20 -> 2016-08-22 02:49:57

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-22 02:49:57"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-22 02:49:53"

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-22 02:49:57"

# The longer "Ready for handoff" label no longer fits the old column
# widths, so re-fit the status/date columns that changed on each sheet
# (matches Excel's own behaviour of widening columns to fit new content).
$wsOverview.Range("E:F").EntireColumn.AutoFit()
$wsZhCn.Range("C:C").EntireColumn.AutoFit()
$wsDeDe.Range("C:C").EntireColumn.AutoFit()
